{"js": "// Update the heading \"6. Technische Beschreibung (Machine Learning Canvas)\"\n// to \"6. Technische Beschreibung (Machine Learning)\" \u2014 i.e. drop \" Canvas\".\n// Target only the \" Canvas\" substring so the rest of the run (rsid\n// attributes, lastRenderedPageBreak, formatting) is left untouched.\nconst results = context.document.body.search(\" Canvas\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n  target.delete();\n  await context.sync();\n}\n", "ps1": "# Update the heading \"6. Technische Beschreibung (Machine Learning Canvas)\"\n# to \"6. Technische Beschreibung (Machine Learning)\" \u2014 i.e. drop \" Canvas\".\n$d = $word.ActiveDocument\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \" Canvas\"\n$rng.Find.MatchCase = $true\n$rng.Find.Execute() | Out-Null\nif ($rng.Find.Found) {\n    $rng.Delete()\n}\n"}
